$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Rows near the end of the table (original indices 34, 35, 36) ---
# Row 34: "20<tab>0.00002<tab>0.00004<tab>0.00003<tab>0.00001<tab>0.00002<tab>0.00003<tab>0.00003<tab>0.00054<tab>100.0"
#   -> collapsed to a single run of "100"
$t.Rows.Item(34).Cells.Item(1).Range.Text = "100"

# Row 35: "1<tab>0.00004<tab>0.00004<tab>0.00004<tab>0.00000<tab>0.00004<tab>0.00004<tab>0.00004<tab>0.00004<tab>100.0"
#   -> collapsed to a single run of "0"
$t.Rows.Item(35).Cells.Item(1).Range.Text = "0"

# Row 36: was an empty run (no <w:t> at all) -> gets text "70"
$t.Rows.Item(36).Cells.Item(1).Range.Text = "70"

# --- Rows near the start of the table (original indices 1, 2, 3) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Insert 10 new rows right after (original) row 3, before (original) row 4 ---
$newValues = @("21", "0.00002", "0.00004", "0.00004", "0.00001", "0.00004", "0.00004", "0.00004", "0.00058", "100.0")
$insertBeforeIndex = 4
foreach ($val in $newValues) {
    $newRow = $t.Rows.Add($t.Rows.Item($insertBeforeIndex))
    $newRow.Cells.Item(1).Range.Text = $val
    $insertBeforeIndex = $insertBeforeIndex + 1
}

Write-Output ("Final row count: " + $t.Rows.Count)
